# Fix formatting issues introduced when scrapping data:
#  1) A handful of "Razon social"/"Nombre Fantasia" entries used a comma where a
#     period was intended as a separator between names.
#  2) The "Importe" column (H) holds amounts that were scraped using Spanish
#     number formatting (thousands separator "." and decimal separator ",").
#     They need to read as plain numeric text using "." as the decimal mark and
#     no thousands separator (e.g. "4.516,00" -> "4516.00").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Razon social / Nombre Fantasia comma -> period fixes
# ---------------------------------------------------------------------------
$textFixes = @{
    "E33"  = "FERNANDEZ MARIO H. GALLICET OSCAR M"
    "E97"  = "FERNANDEZ MARIO H. GALLICET OSCAR M"
    "E98"  = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
    "F98"  = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
    "E136" = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
    "F136" = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
    "E99"  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "E175" = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
    "E115" = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "F115" = "TRABICHET MARIA. VERGARA ADEL Y OTRA"
    "E131" = "RICCOTTI. MARIANA EDITH"
    "E206" = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
}

foreach ($addr in $textFixes.Keys) {
    $ws.Range($addr).Value = $textFixes[$addr]
}

# ---------------------------------------------------------------------------
# 2) Importe column (H2:H238): re-format Spanish-style numeric text as plain
#    numeric text, preserving the original (text) cell type.
# ---------------------------------------------------------------------------
$importeRange = $ws.Range("H2:H238")
# Pre-format as Text so assigning a numeric-looking string keeps it as text
# (otherwise Excel would coerce it into a real number and drop the ".00").
$importeRange.NumberFormat = "@"

for ($r = 2; $r -le 238; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Text
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.Value = $new
}
